$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PQ Challenge 5: rank the B3:B52 "wrapped" country block (5 fields x 10
# countries stored as a flat column) by 2022 population, descending, and
# spill the ranked table (Rank, Country, Population, Growth Rate, Area,
# Density) starting at M2 using a single LAMBDA/LET dynamic-array formula.
$ws.Range("M2").Formula2 = "=LAMBDA(Data, LET(Data2, B3:B52, TransformToCol, WRAPROWS(Data, 5), SortedData, SORTBY(TransformToCol, CHOOSECOLS(TransformToCol, 2), -1), Rank, SEQUENCE(ROWS(SortedData)), Result, HSTACK(Rank, SortedData), Result))(B3:B52)"

# Autofit the newly-populated helper columns (M:O) the way Excel does
# after entering a new block of data.
$ws.Columns.Item(13).AutoFit()
$ws.Columns.Item(14).AutoFit()
$ws.Columns.Item(15).AutoFit()

# Move/restore the active selection to the new formula cell.
$ws.Range("M2").Select()
